$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (HTH / 296001 / GRD / Sheet-2 / 200 / TDFL STRT / 1-A-R / 776)
# This shifts rows 3-6 up to become rows 2-5.
$ws.Rows("2:2").Delete()
